$wb = $excel.ActiveWorkbook

# --- Reorder sheets: original order is GNG, NB, RS, TOL, vSAT ---
# --- Target order: vSAT, TOL, NB, RS, GNG ---
$wsVSAT = $wb.Worksheets.Item("vSAT_TO-16512555236660988")
$wsVSAT.Move($wb.Worksheets.Item(1))
$wsTOL = $wb.Worksheets.Item("TOL_TO-16512555235870974")
$wsTOL.Move($wb.Worksheets.Item(2))
$wsNB = $wb.Worksheets.Item("NB_TO-16512555235230985")
$wsNB.Move($wb.Worksheets.Item(3))
$wsRS = $wb.Worksheets.Item("RS_TO-16512555235311027")
$wsRS.Move($wb.Worksheets.Item(4))

# --- Rename sheets in their new positions & update task-order file lists ---

$ws = $wb.Worksheets.Item(1)
$ws.Name = "vSAT_TO-16515889504398339"
$ws.Range("B2").Value = "SAT_stims-1651588950392177.csv"
$ws.Range("B3").Value = "SAT_stims-1651588950367156.csv"
$ws.Range("B4").Value = "vSAT_stims-16515889504082508.csv"
$ws.Range("B5").Value = "vSAT_stims-1651588950424242.csv"

$ws = $wb.Worksheets.Item(2)
$ws.Name = "TOL_TO-16515889505034037"
$ws.Range("B2").Value = "MM_stims-16515889504715822.csv"
$ws.Range("B3").Value = "ZM_stims-16515889504447474.csv"
$ws.Range("B4").Value = "MM_stims-16515889504864054.csv"
$ws.Range("B5").Value = "ZM_stims-1651588950472588.csv"
$ws.Range("B6").Value = "MM_stims-16515889505024047.csv"
$ws.Range("B7").Value = "ZM_stims-16515889504874077.csv"

$ws = $wb.Worksheets.Item(3)
$ws.Name = "NB_TO-16515889544645734"
$ws.Range("B2").Value = "OB-16515889508287427.csv"
$ws.Range("B3").Value = "TB-16515889544404807.csv"
$ws.Range("B4").Value = "ZB-match_2-16515889505481925.csv"
$ws.Range("B5").Value = "TB-16515889530840511.csv"
$ws.Range("B6").Value = "ZB-match_4-16515889505907612.csv"
$ws.Range("B7").Value = "OB-1651588950870112.csv"
$ws.Range("B8").Value = "TB-16515889535688894.csv"
$ws.Range("B9").Value = "OB-16515889516564379.csv"
$ws.Range("B10").Value = "ZB-match_7-16515889507180872.csv"

$ws = $wb.Worksheets.Item(4)
$ws.Name = "RS_TO-1651588954466575"
$ws.Range("B2").Value = "eyes closed"
$ws.Range("B3").Value = "eyes open"

$ws = $wb.Worksheets.Item(5)
$ws.Name = "GNG_TO-16515889544965906"
$ws.Range("B2").Value = "go_stims-1651588954468603.csv"
$ws.Range("B3").Value = "GNG_stims-1651588954480534.csv"
$ws.Range("B4").Value = "go_stims-1651588954481509.csv"
$ws.Range("B5").Value = "GNG_stims-16515889544965906.csv"
